$wb = $excel.ActiveWorkbook

# The same data-table edits apply to both the "展览" sheet and the "全部类型" sheet,
# which mirror each other in this workbook.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(2, 6).Value = 280   # F2: 279 -> 280
    $ws.Cells.Item(2, 7).Value = 45   # G2: 36.6 -> 45
    $ws.Cells.Item(3, 6).Value = 1429   # F3: 1420 -> 1429
    $ws.Cells.Item(4, 6).Value = 172   # F4: 169 -> 172
    $ws.Cells.Item(7, 6).Value = 104   # F7: 103 -> 104
    $ws.Cells.Item(8, 6).Value = 20   # F8: 19 -> 20
    $ws.Cells.Item(9, 6).Value = 195   # F9: 194 -> 195
    $ws.Cells.Item(11, 6).Value = 3   # F11: 0 -> 3
    $ws.Cells.Item(12, 6).Value = 4744   # F12: 4728 -> 4744
    $ws.Cells.Item(13, 6).Value = 1   # F13: 0 -> 1
    $ws.Cells.Item(14, 6).Value = 7008   # F14: 6987 -> 7008
    $ws.Cells.Item(18, 6).Value = 581   # F18: 579 -> 581
    $ws.Cells.Item(19, 6).Value = 59   # F19: 58 -> 59
    $ws.Cells.Item(21, 6).Value = 4179   # F21: 4177 -> 4179
    $ws.Cells.Item(22, 6).Value = 1114   # F22: 1075 -> 1114
    $ws.Cells.Item(23, 6).Value = 80   # F23: 78 -> 80
    $ws.Cells.Item(24, 6).Value = 77   # F24: 74 -> 77
    $ws.Cells.Item(25, 6).Value = 2757   # F25: 2754 -> 2757
    $ws.Cells.Item(27, 6).Value = 554   # F27: 552 -> 554
    $ws.Cells.Item(29, 6).Value = 386   # F29: 384 -> 386
    $ws.Cells.Item(30, 6).Value = 388   # F30: 386 -> 388
    $ws.Cells.Item(32, 6).Value = 243   # F32: 242 -> 243
    $ws.Cells.Item(34, 6).Value = 1649   # F34: 1648 -> 1649
    $ws.Cells.Item(35, 6).Value = 1060   # F35: 1055 -> 1060
    $ws.Cells.Item(37, 6).Value = 725   # F37: 691 -> 725
    $ws.Cells.Item(39, 6).Value = 557   # F39: 554 -> 557
    $ws.Cells.Item(41, 6).Value = 499   # F41: 498 -> 499
    $ws.Cells.Item(43, 6).Value = 18   # F43: 10 -> 18
    $ws.Cells.Item(45, 6).Value = 255   # F45: 238 -> 255
    $ws.Cells.Item(46, 6).Value = 658   # F46: 656 -> 658
    $ws.Cells.Item(47, 6).Value = 23   # F47: 22 -> 23

    # I45 cover image URL update
    $ws.Cells.Item(45, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/O2VuWJFW1721291227478.jpeg"
}

